# Refresh the cryptos price/volume snapshot (scheduled GitHub Actions update).
# Source data is plain text ("Price"/"Volume(1h)" columns are formatted strings,
# not numbers), so every write below targets a Text-typed cell. Numeric-looking
# price strings get a leading apostrophe so Excel keeps them as text (otherwise
# trailing zeros such as "0.200" or "5.30" would be silently dropped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '69.864.00'  # D (Price)
$ws.Cells.Item(2, 5).Value = '  +0.79%  '  # E (Volume 1h)

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '3.524.96'  # D (Price)

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).Value = '''0.998'  # D (Price)
$ws.Cells.Item(4, 5).Value = '  -0.04%  '  # E (Volume 1h)

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = '''606.51'  # D (Price)
$ws.Cells.Item(5, 5).Value = '  -0.63%  '  # E (Volume 1h)

# Row 6: Solana
$ws.Cells.Item(6, 4).Value = '''197.41'  # D (Price)
$ws.Cells.Item(6, 5).Value = '  +6.42%  '  # E (Volume 1h)

# Row 7: XRP
$ws.Cells.Item(7, 5).Value = '  -0.58%  '  # E (Volume 1h)

# Row 8: USDC
$ws.Cells.Item(8, 4).Value = '''0.999'  # D (Price)
$ws.Cells.Item(8, 5).Value = '  -0.06%  '  # E (Volume 1h)

# Row 9: Dogecoin
$ws.Cells.Item(9, 4).Value = '''0.200'  # D (Price)
$ws.Cells.Item(9, 5).Value = '  -7.57%  '  # E (Volume 1h)

# Row 10: Cardano
$ws.Cells.Item(10, 4).Value = '''0.649'  # D (Price)
$ws.Cells.Item(10, 5).Value = '  -0.49%  '  # E (Volume 1h)

# Row 11: Avalanche
$ws.Cells.Item(11, 4).Value = '''53.83'  # D (Price)
$ws.Cells.Item(11, 5).Value = '  +1.63%  '  # E (Volume 1h)

# Row 12: ShibaInu
$ws.Cells.Item(12, 5).Value = '  -2.29%  '  # E (Volume 1h)

# Row 14: WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = '4.083.97'  # D (Price)
$ws.Cells.Item(14, 5).Value = '  +1.08%  '  # E (Volume 1h)

# Row 15: BitcoinCash
$ws.Cells.Item(15, 4).Value = '''599.22'  # D (Price)
$ws.Cells.Item(15, 5).Value = '  -0.35%  '  # E (Volume 1h)

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = '70.024.42'  # D (Price)
$ws.Cells.Item(16, 5).Value = '  +0.98%  '  # E (Volume 1h)

# Row 17: Chainlink
$ws.Cells.Item(17, 4).Value = '''19.11'  # D (Price)
$ws.Cells.Item(17, 5).Value = '  +1.51%  '  # E (Volume 1h)

# Row 18: Uniswap
$ws.Cells.Item(18, 4).Value = '''12.77'  # D (Price)
$ws.Cells.Item(18, 5).Value = '  +1.45%  '  # E (Volume 1h)

# Row 19: WrappedEther
$ws.Cells.Item(19, 4).Value = '3.532.97'  # D (Price)
$ws.Cells.Item(19, 5).Value = '  +1.04%  '  # E (Volume 1h)

# Row 20: TRON
$ws.Cells.Item(20, 5).Value = '  +1.34%  '  # E (Volume 1h)

# Row 21: Polygon
$ws.Cells.Item(21, 4).Value = '''0.994'  # D (Price)
$ws.Cells.Item(21, 5).Value = '  +0.75%  '  # E (Volume 1h)

# Row 22: InternetComputer(DFINITY)
$ws.Cells.Item(22, 5).Value = '  +7.23%  '  # E (Volume 1h)

# Row 23: Toncoin
$ws.Cells.Item(23, 4).Value = '''5.30'  # D (Price)
$ws.Cells.Item(23, 5).Value = '  +4.78%  '  # E (Volume 1h)

# Row 24: Litecoin
$ws.Cells.Item(24, 4).Value = '''101.99'  # D (Price)
$ws.Cells.Item(24, 5).Value = '  -3.47%  '  # E (Volume 1h)

# Row 25: PancakeSwap
$ws.Cells.Item(25, 4).Value = '''4.64'  # D (Price)
$ws.Cells.Item(25, 5).Value = '  +0.01%  '  # E (Volume 1h)

# Row 26: ImmutableX
$ws.Cells.Item(26, 4).Value = '''3.19'  # D (Price)
$ws.Cells.Item(26, 5).Value = '  +5.92%  '  # E (Volume 1h)

# Row 27: RenderToken
$ws.Cells.Item(27, 5).Value = '  -0.56%  '  # E (Volume 1h)

# Row 28: Filecoin
$ws.Cells.Item(28, 5).Value = '  -1.91%  '  # E (Volume 1h)

# Row 29: EthereumClassic
$ws.Cells.Item(29, 4).Value = '''33.39'  # D (Price)
$ws.Cells.Item(29, 5).Value = '  -0.46%  '  # E (Volume 1h)

# Row 30: dogwifhat
$ws.Cells.Item(30, 4).Value = '''4.33'  # D (Price)
$ws.Cells.Item(30, 5).Value = '  +10.95%  '  # E (Volume 1h)

# Row 31: NEARProtocol
$ws.Cells.Item(31, 4).Value = '''7.08'  # D (Price)
$ws.Cells.Item(31, 5).Value = '  +1.82%  '  # E (Volume 1h)

# Row 32: Cosmos
$ws.Cells.Item(32, 5).Value = '  +0.89%  '  # E (Volume 1h)

# Row 33: Hedera
$ws.Cells.Item(33, 5).Value = '  -0.35%  '  # E (Volume 1h)

# Row 34: OKB
$ws.Cells.Item(34, 4).Value = '''63.14'  # D (Price)
$ws.Cells.Item(34, 5).Value = '  -0.10%  '  # E (Volume 1h)

# Row 35: PEPE
$ws.Cells.Item(35, 4).Value = '0.0₃0862'  # D (Price)
$ws.Cells.Item(35, 5).Value = '  +10.95%  '  # E (Volume 1h)

# Row 36: Maker
$ws.Cells.Item(36, 4).Value = '3.728.83'  # D (Price)
$ws.Cells.Item(36, 5).Value = '  +3.71%  '  # E (Volume 1h)

# Row 37: Dai -> Fetch.AI
$ws.Cells.Item(37, 2).Value = 'Fetch.AI'  # B (Coin)
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'  # C (Link)
$ws.Cells.Item(37, 4).Value = '''3.08'  # D (Price)
$ws.Cells.Item(37, 5).Value = '  -3.48%  '  # E (Volume 1h)

# Row 38: Fetch.AI -> Dai
$ws.Cells.Item(38, 2).Value = 'Dai'  # B (Coin)
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'  # C (Link)
$ws.Cells.Item(38, 4).Value = '''1.00'  # D (Price)
$ws.Cells.Item(38, 5).Value = '  +0.32%  '  # E (Volume 1h)

# Row 39: Stacks
$ws.Cells.Item(39, 5).Value = '  -0.01%  '  # E (Volume 1h)

# Row 40: TheGraph
$ws.Cells.Item(40, 4).Value = '''0.394'  # D (Price)
$ws.Cells.Item(40, 5).Value = '  -0.77%  '  # E (Volume 1h)

# Row 41: InjectiveProtocol
$ws.Cells.Item(41, 4).Value = '''36.68'  # D (Price)
$ws.Cells.Item(41, 5).Value = '  -0.08%  '  # E (Volume 1h)

# Row 42: Bittensor
$ws.Cells.Item(42, 4).Value = '''488.08'  # D (Price)
$ws.Cells.Item(42, 5).Value = '  -5.77%  '  # E (Volume 1h)

# Row 43: Kaspa
$ws.Cells.Item(43, 5).Value = '  -3.59%  '  # E (Volume 1h)

# Row 44: VeChain
$ws.Cells.Item(44, 5).Value = '  -1.18%  '  # E (Volume 1h)

# Row 45: ThetaToken
$ws.Cells.Item(45, 5).Value = '  -3.25%  '  # E (Volume 1h)

# Row 47: ApeXProtocol
$ws.Cells.Item(47, 5).Value = '  -1.26%  '  # E (Volume 1h)

# Row 48: FirstDigitalUSD
$ws.Cells.Item(48, 5).Value = '  +0.35%  '  # E (Volume 1h)

# Row 49: THORChain
$ws.Cells.Item(49, 5).Value = '  -3.16%  '  # E (Volume 1h)

# Row 50: FLOKI
$ws.Cells.Item(50, 5).Value = '  +2.32%  '  # E (Volume 1h)

# Row 51: Monero -> Mantle
$ws.Cells.Item(51, 2).Value = 'Mantle'  # B (Coin)
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'  # C (Link)
$ws.Cells.Item(51, 4).Value = '''1.29'  # D (Price)
$ws.Cells.Item(51, 5).Value = '  +11.22%  '  # E (Volume 1h)
